# Ajustes de los botones de cada universidad: se agrega la columna "universidades"
# (columna U) con el listado de universidades asociadas a las primeras ofertas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nueva columna U: encabezado + datos
$ws.Range("U1").Value = "universidades"
$ws.Range("U2").Value = "Areandina"
$ws.Range("U3").Value = "militar"
$ws.Range("U4").Value = "Simonbolivar"

# Copiar el formato del encabezado (negrita + bordes) de T1 a U1
$ws.Range("T1").Copy() | Out-Null
$ws.Range("U1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Mover la selección a la celda recién editada, como quedó en el archivo original
$ws.Range("U4").Select() | Out-Null
